$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at position 14 (N), shifting old N,O to O,P
$ws.Columns.Item(14).Insert()

# New column N header
$ws.Range("N1").Value = "subject"

# Comments are anchored to cell position, not content, so they did not
# shift together with the cell values during the column insert. Capture
# the two affected comments' text, delete them, then re-create them one
# column to the right (N1 -> O1, O1 -> P1).
$txtN1 = $ws.Range("N1").Comment.Text()
$txtO1 = $ws.Range("O1").Comment.Text()

$ws.Range("N1").Comment.Delete()
$ws.Range("O1").Comment.Delete()

$ws.Range("O1").AddComment($txtN1) | Out-Null
$ws.Range("P1").AddComment($txtO1) | Out-Null

# New column N gets the same (non-bestFit) width as column M
$ws.Columns.Item(14).ColumnWidth = 13

# New hyperlink on N1 (property-subject)
$ws.Hyperlinks.Add($ws.Range("N1"), "https://docs.dasch.swiss/latest/DSP-TOOLS/file-formats/json-project/ontologies/", "property-subject") | Out-Null

# Update active selection
$ws.Range("D32").Select() | Out-Null

Write-Output "column inserted"
